$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.702.74'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.49%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.455.10'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.61%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '571.12'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.08'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.39%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.528'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.80%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.84%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.15'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.27%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.29%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '28.64'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.12%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000173'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.82%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.895.57'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '62.596.14'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.77%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.456.59'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.53%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.68'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -6.15%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.70'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -2.86%  '
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = 'SuiNetwork'
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.23'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.13%  '
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '320.99'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.57%  '
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = 'Polkadot'
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.13'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.92'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +3.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '64.73'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.25%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '648.19'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.55%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.579.57'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.45%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.23%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0947'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -3.89%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.40'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.74%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.79'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -3.45%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.80'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.16%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.45%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.06%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.48'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.63%  '
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.61'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -3.48%  '
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '151.79'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.54%  '
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = 'PolygonEcosystemToken'
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.363'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -2.21%  '
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = 'EthereumClassic'
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.48'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.44%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.28'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.49%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.62'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -3.90%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.69'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -3.62%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0₆0306'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.49%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '152.25'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.58%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '15.40'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.73%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.52'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.06%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.602'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.50%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '19.85'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -3.87%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0501'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.26%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0900'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.94%  '
